$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Recorded By" (column G) values for data rows 2-28 and 30
# (row 29 is already empty in the source data).
$rowsToClear = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,30)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 7).Value = ""
}

# Shrink column G width from 50 to 13 (stored width = ColumnWidth + 0.83,
# so 12.17 round-trips to the target stored width of 13)
$ws.Columns.Item(7).ColumnWidth = 12.17
